$d = $word.ActiveDocument

function Replace-AllText($find, $replace) {
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $replace, 2) | Out-Null
}

# Title heading (appears twice: H1 title and bold summary line at the end)
Replace-AllText "Play Magic Oak Slot Game for Free - Exciting Bonus Features" `
                 "Play Magic Oak Free - Exciting Gameplay and Generous Payouts"

# "What we like" bullet list
Replace-AllText "Great graphics and animation" `
                 "Variety of betting options and limits"

Replace-AllText "Exciting bonus features with Wild and Scatter symbols" `
                 "Exciting gameplay mechanics and features"

Replace-AllText "Auto-spin feature for up to 500 spins" `
                 "Generous payouts for winning combinations"

Replace-AllText "High RTP of 96.71%" `
                 "Free spins feature adds extra excitement"

# "What we don't like" bullet list
Replace-AllText "Limited betting options for high rollers" `
                 "Limited number of paylines"

Replace-AllText "Fewer paylines compared to other slot games" `
                 "No progressive jackpot feature"

# Closing italic summary paragraph
Replace-AllText "Get ready to play Magic Oak - a 4x4 slot game with 20 fixed paylines, Wild and Scatter symbols, and exciting bonus features. Play for free now!" `
                 "Play Magic Oak for free and enjoy exciting gameplay mechanics and generous payouts."
